# CIERRE 4 JUN 22
# Update the "ARQUITECTO" sheet (first sheet) with the new incentive amount
# and the new concept text, and move the active selection to E11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Amount paid: 50,000 -> 150,000
$ws.Range("D1").Value = 150000

# Amount in words: "CINCUENTA MIL..." -> "CIENTO CINCUENTA MIL..."
$ws.Range("A2").Value = "CIENTO CINCUENTA    MIL   PESOS 00/100 M.N."

# Concept: demolition payment -> parking lot construction payment
$ws.Range("A4").Value = "OBRA ESTACIONAMIENTO OBRADOR "

# Move selection to E11 on this (active) sheet
$ws.Activate()
$ws.Range("E11").Select()
